# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2:F3, zh-cn!C2:C3,
#   de-de!C2:C3 all share the same string).
# - Narrow the "Latest HO Xliff Generate Date" columns (Overview E:F,
#   zh-cn C, de-de C) from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text change: "Ready for handoff" -> "In Translation" ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column width changes (narrower date columns) ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
